$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 322.33334
$ws.Range("J19").Value = 94
$ws.Range("L19").Value = 94
$ws.Range("N19").Value = -444
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").Value = ""
$ws.Range("H74").Value = 2916.6667
$ws.Range("I74").Value = 2916.6667
$ws.Range("K74").Value = 2916.6667
$ws.Range("M74").Value = -1980.6667
$ws.Range("H77").Value = 2916.6667
$ws.Range("I77").Value = 2916.6667
$ws.Range("K77").Value = 14583.3335
$ws.Range("M77").Value = -9903.333500000001
$ws.Range("H137").Value = 2334.7856
$ws.Range("I137").Value = 1476.6428
$ws.Range("K137").Value = 4429.928400000001
$ws.Range("M137").Value = -1879.928400000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 1250.5
$ws.Range("J30").Value = 2492
$ws.Range("L30").Value = 2492
$ws.Range("N30").Value = -2792
$ws.Range("H32").Value = 4688.2856
$ws.Range("I32").Value = 4038.926
$ws.Range("K32").Value = 4038.926
$ws.Range("M32").Value = -3751.926
$ws.Range("H61").Value = 8999.166999999999
$ws.Range("I61").Value = 10999.333
$ws.Range("K61").Value = 10999.333
$ws.Range("M61").Value = -10787.333
$ws.Range("H132").Value = 4337.25
$ws.Range("I132").Value = 3939.8
$ws.Range("K132").Value = 11819.4
$ws.Range("M132").Value = -9289.400000000001
$ws.Range("H136").Value = 8999.166999999999
$ws.Range("I136").Value = 10999.333
$ws.Range("K136").Value = 32997.999
$ws.Range("M136").Value = -30447.999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 3800490
$ws.Range("I7").Value = 9500500
$ws.Range("K7").Value = 9500500
$ws.Range("M7").Value = -9500387
$ws.Range("H20").Value = 2645.2666
$ws.Range("I20").Value = 1447.2727
$ws.Range("K20").Value = 1447.2727
$ws.Range("M20").Value = -1200.2727
$ws.Range("H94").Value = 875.8570999999999
$ws.Range("I94").Value = 875.8570999999999
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 875.8570999999999
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -424.8570999999999
$ws.Range("N94").Value = ""
$ws.Range("H99").Value = 2210.3635
$ws.Range("I99").Value = 2123.4285
$ws.Range("K99").Value = 2123.4285
$ws.Range("M99").Value = -625.4285
$ws.Range("H105").Value = 5025.9
$ws.Range("I105").Value = 3752.5
$ws.Range("K105").Value = 3752.5
$ws.Range("M105").Value = -2005.5
$ws.Range("H107").Value = 1205.8667
$ws.Range("I107").Value = 1262.7858
$ws.Range("K107").Value = 1262.7858
$ws.Range("M107").Value = 657.2141999999999
$ws.Range("H134").Value = 2611.375
$ws.Range("I134").Value = 1829.3846
$ws.Range("K134").Value = 5488.1538
$ws.Range("M134").Value = -2953.1538
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 1200
$ws.Range("I21").Value = 1200
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 1200
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -965
$ws.Range("N21").Value = ""
$ws.Range("H31").Value = 5886.6
$ws.Range("I31").Value = 3809
$ws.Range("J31").Value = 9003
$ws.Range("K31").Value = 3809
$ws.Range("L31").Value = 9003
$ws.Range("M31").Value = -3514
$ws.Range("N31").Value = -9593
$ws.Range("H34").Value = 5886.6
$ws.Range("I34").Value = 3809
$ws.Range("J34").Value = 9003
$ws.Range("K34").Value = 3809
$ws.Range("L34").Value = 9003
$ws.Range("M34").Value = -3607
$ws.Range("N34").Value = -9407
$ws.Range("H105").Value = 2250
$ws.Range("I105").Value = 1500
$ws.Range("K105").Value = 1500
$ws.Range("M105").Value = 247
$ws.Range("H132").Value = 7563.5625
$ws.Range("I132").Value = 5961.4
$ws.Range("K132").Value = 17884.2
$ws.Range("M132").Value = -15354.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 901.2857
$ws.Range("I34").Value = 40
$ws.Range("J34").Value = 1044.8334
$ws.Range("K34").Value = 120
$ws.Range("L34").Value = 3134.5002
$ws.Range("M34").Value = -36
$ws.Range("N34").Value = -3302.5002
$ws.Range("H104").Value = 14999
$ws.Range("J104").Value = 14999
$ws.Range("L104").Value = 44997
$ws.Range("N104").Value = -50239
$ws.Range("H113").Value = 371
$ws.Range("J113").Value = 371
$ws.Range("L113").Value = 1113
$ws.Range("N113").Value = -5453
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 2000
$ws.Range("J12").Value = 2000
$ws.Range("L12").Value = 2000
$ws.Range("N12").Value = -2280
$ws.Range("H22").Value = 1710.3334
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = ""
$ws.Range("H70").Value = 5795.8
$ws.Range("I70").Value = 5491
$ws.Range("K70").Value = 5491
$ws.Range("M70").Value = -5221
$ws.Range("H73").Value = 5795.8
$ws.Range("I73").Value = 5491
$ws.Range("K73").Value = 5491
$ws.Range("M73").Value = -4555
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1042.2
$ws.Range("I16").Value = 1042.2
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1042.2
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -872.2
$ws.Range("N16").Value = ""
$ws.Range("H22").Value = 1234.6666
$ws.Range("I22").Value = 1209.5
$ws.Range("J22").Value = 1285
$ws.Range("K22").Value = 1209.5
$ws.Range("L22").Value = 1285
$ws.Range("M22").Value = -914.5
$ws.Range("N22").Value = -1875
$ws.Range("H27").Value = 1234.6666
$ws.Range("I27").Value = 1209.5
$ws.Range("J27").Value = 1285
$ws.Range("K27").Value = 1209.5
$ws.Range("L27").Value = 1285
$ws.Range("M27").Value = -1102.5
$ws.Range("N27").Value = -1499
$ws.Range("H40").Value = 2748.5
$ws.Range("I40").Value = 2748.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2748.5
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2612.5
$ws.Range("N40").Value = ""
$ws.Range("H46").Value = 1739.8387
$ws.Range("J46").Value = 2012.579
$ws.Range("L46").Value = 2012.579
$ws.Range("N46").Value = -2388.579
$ws.Range("H140").Value = 73607.25
$ws.Range("J140").Value = 73607.25
$ws.Range("L140").Value = 73607.25
$ws.Range("N140").Value = -83967.25
$ws.Range("H141").Value = 59999
$ws.Range("J141").Value = 59999
$ws.Range("L141").Value = 59999
$ws.Range("N141").Value = -70359
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = ""
$ws.Range("H107").Value = 522.5
$ws.Range("I107").Value = 522.5
$ws.Range("K107").Value = 1567.5
$ws.Range("M107").Value = 352.5

Write-Output "Applied all Maduin_Profits updates"